# Update cryptos list snapshot values (prices / 1h volume change %)
# and correct the ordering/values for three rows (49-51) per the
# upstream GitHub Actions refresh of cryptos.xlsx.
#
# NOTE: Price values are stored as plain text (e.g. "606.50",
# "66.159.08"). Any value that Excel would otherwise auto-parse as a
# number is given an explicit "@" (text) number format *before* the
# value is assigned, so it stays text instead of being coerced into a
# floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.159.08"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "3.567.31"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.50"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.32"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").Value = "3.566.60"
$ws.Range("E7").Value = "  +1.92%  "
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.93"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.413"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "4.171.50"
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000208"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.08"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "3.567.64"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").Value = "66.255.10"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.42"
$ws.Range("E19").Value = "  +8.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.87"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "430.29"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.613"
$ws.Range("E23").Value = "  +3.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.16"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "3.708.56"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.51"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.97"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.14"
$ws.Range("E30").Value = "  -3.90%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.66"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("D34").Value = "3.562.31"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("E35").Value = "  -6.07%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.88"
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.62"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "175.67"
$ws.Range("E41").Value = "  +3.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0852"
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.22"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.894"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.95"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.06"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.84"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.37"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.65"
$ws.Range("E50").Value = "  +7.84%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.14"
$ws.Range("E51").Value = "  -0.10%  "
